$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update wording: "dead time(/dead-time)" -> "nonlinearity" related terms
$ws.Range("B4").Value = "analyze & investigate the nonlineariy effect on the VSIs for various PMSMs"
$ws.Range("B6").Value = "propose a software based nonlinearity compensation method "

# Move the active selection to B6, matching the post-edit cursor position
$ws.Range("B6").Select()
